# TC_15 DemoWebshop_TotalOrders: Written TotalNumberOf Orders and SumOfAllOrder
# values back to excel file.
#
# This script reproduces the authoring edit:
#   - Adds two new result columns (F: TotalNumberOfOrders, G: SumOfAllOrders)
#     to the "DemoWebshop_TotalOrders" sheet, with the header in row 1 and
#     the captured run's results in row 2.
#   - Updates the order number captured on the "DemoWebshop_ReOrder" sheet
#     (a side effect of the same automation run writing its results back).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. DemoWebshop_TotalOrders: add TotalNumberOfOrders / SumOfAllOrders columns
# ---------------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("DemoWebshop_TotalOrders")
$wsTotal.Activate()

# Headers (row 1)
$hdr1 = $wsTotal.Cells.Item(1, 6)
$hdr1.Value = "TotalNumberOfOrders"
$hdr1.Style = "Normal"

$hdr2 = $wsTotal.Cells.Item(1, 7)
$hdr2.Value = "SumOfAllOrders"
$hdr2.Style = "Normal"

# Captured results (row 2) - written back as plain text, matching the
# automation framework's reporting format.
$totalOrders = $wsTotal.Cells.Item(2, 6)
$totalOrders.NumberFormat = "@"
$totalOrders.Value = "313"
$totalOrders.Style = "Normal"

$sumOrders = $wsTotal.Cells.Item(2, 7)
$sumOrders.NumberFormat = "@"
$sumOrders.Value = "47639.6"
$sumOrders.Style = "Normal"

# Column widths for the two new columns.
$wsTotal.Columns.Item(6).ColumnWidth = 22
$wsTotal.Columns.Item(7).ColumnWidth = 19.5

# View state: scroll the grid over and land the selection on the new
# SumOfAllOrders header cell (G1), as happens after typing the last value.
$wsTotal.Range("G1").Select()
$excel.ActiveWindow.ScrollColumn = 3

# Page setup was touched (printing orientation) as part of the save.
$wsTotal.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# 2. DemoWebshop_ReOrder: refresh the captured order number
# ---------------------------------------------------------------------------
$wsReOrder = $wb.Worksheets.Item("DemoWebshop_ReOrder")
$orderCell = $wsReOrder.Cells.Item(2, 6)
$orderCell.NumberFormat = "@"
$orderCell.Value = "Order number: 1481392"
$orderCell.Style = "Normal"

# ---------------------------------------------------------------------------
# 3. Workbook view: bring the DemoWebshop sheets into view in the tab strip
# ---------------------------------------------------------------------------
$wb.Windows.Item(1).ScrollWorkbookTabs(3, 2) | Out-Null

$wsTotal.Activate()
